$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (B1:P1) with full party descriptions ---
$ws.Range("B1").Value = "ASZ - Agrarian Alliance (Agrárszövetség, ASZ)"
$ws.Range("C1").Value = "FKGP - Independent Smallholders Party (Független Kisgazdapárt, FKGP)"
$ws.Range("D1").Value = "Fidesz-MPP - Alliance of Young Democrats-Hungarian Civic Party (Fiatal Demokraták Szövetsége-Magyar Polgári Párt, Fidesz-MPP)"
$ws.Range("E1").Value = "KDNP - Christian-Democratic People's Party (Keresztenydemokrata`n Neppitrt (KDNP), KDNP)"
$ws.Range("F1").Value = "LA - Liberal Alliance - Party of Entrepreneurs (Liberális Szövetség-Vállalkozók Pártja, LA)"
$ws.Range("G1").Value = "MDF - Hungarian Democratic Forum (Magyar Demokrata Fórum, MDF)"
$ws.Range("H1").Value = "MSZP - Hungarian Socialist Party (Magyar Szocialista Pirt, MSZP)"
$ws.Range("I1").Value = "SZDSZ - Alliance of Free Democrats (Szabad Demokrathk Szovetstge, SZDSZ)"
$ws.Range("J1").Value = "MIÉP - Party of Hungarian Justice and Life (Magyar Igazságés Élet Párt, MIÉP)"
$ws.Range("K1").Value = "LMP - Politics Can be Different (Lehet Más a Politika, LMP)"
$ws.Range("L1").Value = "MIÉP-Jobbik - Party of Hungarian Justice and Life-Jobbik The Third Way (Magyar Igazság és Élet Pártja/Jobbik a Harmadik Út, MIÉP-Jobbik)"
$ws.Range("M1").Value = "Other - Other"
$ws.Range("N1").Value = "DK - Democratic Coalition (Demokratikus Koalíció, DK)"
$ws.Range("O1").Value = "Együtt - Together - Party for a New Era (Együtt, Együtt)"
$ws.Range("P1").Value = "MNOO - National Self-Government of Germans in Hungary (Magyarországi Németek Országos Önkormányzata, MNOO)"

# --- New column Q1, styled like the rest of the header row ---
$ws.Range("Q1").Value = "No acronym - Independent (Ind., No acronym)"
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Update data grid for rows 6-8 to reflect the new MIÉP-Jobbik column (L) ---
# Row 6 (2010): clear J6, shift K6/L6 stay, add M6
$ws.Range("J6").ClearContents()
$ws.Range("M6").Value = 0

# Row 7 (2014): clear J7, add L7
$ws.Range("J7").ClearContents()
$ws.Range("L7").Value = 0

# Row 8 (2018): clear J8, add L8, clear M8, add Q8
$ws.Range("J8").ClearContents()
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("Q8").Value = 0
